$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text in A1
$ws.Range("A1").Value = "Sales for the last year (365 days)"

# Update row 3 sales figures to reflect a full year's data
$ws.Range("A3").Value = 100.4
$ws.Range("B3").Value = 1520.91
$ws.Range("C3").Value = 33
$ws.Range("D3").Value = 9.9499999999999993
$ws.Range("F3").Value = 24.95
$ws.Range("G3").Value = 492.75
$ws.Range("K3").Value = 7.9
$ws.Range("L3").Value = 48.9
